# The deck's live theme ("Integral" / "Red Violet" colour scheme) is swapped
# back to the built-in "Office Theme" colour scheme. Re-point every themed
# colour slot (Dark1/Light1/Dark2/Light2/Accent1-6/Hyperlink/FollowedHyperlink)
# at the stock Office palette via the slide master's theme colour scheme.

$pres = $ppt.ActivePresentation
$master = $pres.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0        # Dark1    -> 000000
$colors.Item(2).RGB  = 16777215 # Light1   -> FFFFFF
$colors.Item(3).RGB  = 6968388  # Dark2    -> 44546A
$colors.Item(4).RGB  = 15132391 # Light2   -> E7E6E6
$colors.Item(5).RGB  = 13998939 # Accent1  -> 5B9BD5
$colors.Item(6).RGB  = 3243501  # Accent2  -> ED7D31
$colors.Item(7).RGB  = 10855845 # Accent3  -> A5A5A5
$colors.Item(8).RGB  = 49407    # Accent4  -> FFC000
$colors.Item(9).RGB  = 12874308 # Accent5  -> 4472C4
$colors.Item(10).RGB = 4697456  # Accent6  -> 70AD47
$colors.Item(11).RGB = 12673797 # Hyperlink -> 0563C1
$colors.Item(12).RGB = 7491477  # FollowedHyperlink -> 954F72
